$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Description" header in column F, row 1.
$ws.Range("F1").Value = "Description"

# Match the header formatting used by the other header cells in row 1
# (A1/E1 use the "accent" style) by copying A1's format onto F1.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Widen the new column F.
$ws.Columns("F").ColumnWidth = 24.65

# Update the saved selection/active cell for the sheet.
$ws.Range("F4").Select() | Out-Null
